$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# --- Part replacement: diode SS34 -> B340A (row 15) ---
$ws.Range("C15").Value = "B340A"

# --- Quantity / price corrections that ripple into the H column formulas ---
$ws.Range("F11").Value = 5
$ws.Range("F15").Value = 10
$ws.Range("G15").Value = 0.985
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 1.15
$ws.Range("F19").Value = 3

# --- New BOM row 23: fuse (bezpiecznik) B340A fast 2.5A ---
$ws.Range("B23").Value = "bezpiecznik"
$ws.Range("C23").Value = "fast 2.5A"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 1.46
$ws.Range("H23").Formula = "=F23*G23"
$ws.Range("I23").Value = "TME"
$ws.Range("J23").Value = "mouser"

# Hyperlinks for the new fuse row (TME page first, keep its real URL as the
# display text, then restore the short "TME" label in the cell itself)
$ws.Hyperlinks.Add($ws.Range("I23"), "https://www.tme.eu/pl/details/erbrg2r50v/bezpieczniki-smd-1206-szybkie/panasonic/", "", "", "https://www.tme.eu/pl/details/erbrg2r50v/bezpieczniki-smd-1206-szybkie/panasonic/") | Out-Null
$ws.Range("I23").Value = "TME"

$ws.Hyperlinks.Add($ws.Range("J23"), "https://www.mouser.pl/ProductDetail/Bel-Fuse/0ERB-R250-A?qs=fVJ5M%252Bpe2yVH1zzH%2Fc6bsg%3D%3D") | Out-Null
$ws.Range("J23").Value = "mouser"

# Match formatting of the existing rows by copying formats across (done last
# so the plain "Hyperlink" style index used elsewhere in the sheet wins over
# the one Hyperlinks.Add auto-creates)
$ws.Range("C17:D17").Copy()
$ws.Range("C23:D23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23:D23").Merge()

$ws.Range("I19").Copy()
$ws.Range("I23").PasteSpecial(-4122) | Out-Null

$ws.Range("J17").Copy()
$ws.Range("J23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Grand total row moves down to row 24 ---
$ws.Range("H24").Formula = "=SUM(H4:H23)"

# --- Selection bookkeeping that Excel records on save ---
$ws.Activate()
$ws.Range("L13").Select()
